# "update christi data binding"
#
# - sheet1 ("change password"): several "show pwd*" cells flip from
#   "yes" to "no", and a new test-case row (row 8) is appended using a
#   new shared string "Password2#".
# - the active tab moves from sheet3 ("edit reward") back to sheet1
#   ("change password"), with the selection anchored at A8 (the new row).

$wb = $excel.ActiveWorkbook

$wsChangePassword = $wb.Worksheets.Item("change password")

# Flip these "show pwd" / "show pwd baru" / "show pwd konfir" flags from
# yes -> no.
$wsChangePassword.Range("E2").Value = "no"
$wsChangePassword.Range("F3").Value = "no"
$wsChangePassword.Range("G4").Value = "no"
$wsChangePassword.Range("E5").Value = "no"
$wsChangePassword.Range("F5").Value = "no"
$wsChangePassword.Range("F6").Value = "no"
$wsChangePassword.Range("G6").Value = "no"
$wsChangePassword.Range("E7").Value = "no"
$wsChangePassword.Range("G7").Value = "no"

# Append the new test case in row 8.
$wsChangePassword.Range("B8").Value = "Password1#"
$wsChangePassword.Range("C8").Value = "Password2#"
$wsChangePassword.Range("D8").Value = "Password2#"
$wsChangePassword.Range("E8").Value = "yes"
$wsChangePassword.Range("F8").Value = "yes"
$wsChangePassword.Range("G8").Value = "yes"
$wsChangePassword.Range("H8").Value = "pass"

# Move the active tab/selection to "change password" / A8, which takes it
# away from "edit reward".
$wsChangePassword.Activate()
$wsChangePassword.Range("A8").Select()
